# Minor adjustments; Added a country-specific parameter SAVINGS_RATE to be
# loaded from the parameters.xlsx input file.
#
# This script adds a new "Info" worksheet (after the existing
# "ColumnsNumberParameters" sheet), populates it with two descriptive rich
# text notes, and makes it the active/selected sheet - mirroring the target
# OOXML diff.

$wb = $excel.ActiveWorkbook
$wsParams = $wb.Worksheets.Item(1)

# Add the new "Info" sheet right after the parameters sheet.
$wsInfo = $wb.Worksheets.Add($null, $wsParams)
$wsInfo.Name = "Info"

# --- Row 1: explanation of the file's purpose -------------------------
$text1 = "This Excel file is used to define the column numbers required for the corresponding processes in the reg_estimates files."
$wsInfo.Range("A1").Value = $text1

$italicStart1 = $text1.IndexOf("reg_estimates") + 1
$italicLen1 = "reg_estimates".Length
$run1 = $wsInfo.Range("A1").Characters($italicStart1, $italicLen1)
$run1.Font.Italic = $true

# --- Row 2: note about per-country configuration -----------------------
$text2 = "Since the column numbers differ by country, they must be configured separately for each country."
$wsInfo.Range("A2").Value = $text2

$sepStart = $text2.IndexOf("separately") + 1
$sepLen = "separately".Length
$runSep = $wsInfo.Range("A2").Characters($sepStart, $sepLen)
$runSep.Font.Italic = $true
$runSep.Font.Underline = $true

$eachStart = $text2.IndexOf("each country") + 1
$eachLen = "each country".Length
$runEach = $wsInfo.Range("A2").Characters($eachStart, $eachLen)
$runEach.Font.Italic = $true
$runEach.Font.Underline = $true

# Selection inside the Info sheet, and make it the active sheet/tab.
[void]$wsInfo.Range("B6").Select()
[void]$wsInfo.Activate()
